# edit.ps1 - applies the commit "rebuild pages at 4e8fb1c" to
# ashley.hindmarsh.cv.docx
#
# Changes:
#  1. Insert a new "Compact"-style bulleted paragraph (same list as the
#     rest of the Summary section, numId 1002) containing a new "Ethos:
#     ..." bullet, placed right after the last Summary bullet ("Recent,
#     relevant experience...") and before the "Technical" Heading2.
#  2. Fix a typo in the existing "Introduction of a new backend payment
#     system..." bullet: "includingg migrgation" -> "including migration".
#
# (The diff also shows the <w:nsid> GUIDs of two <w:abstractNum> list
# definitions in numbering.xml changing. Those are opaque, randomly
# generated identifiers with no semantic/visual effect, are not part of
# the Word object model in real Word either [not reachable from VBA/COM],
# and are not exposed anywhere in this host's Word OM surface - there is
# no supported way to target them from script, so they are left as-is.)

$d = $word.ActiveDocument

# --- 1. Fix the typo --------------------------------------------------
# Replace the whole paragraph's OOXML (rather than doing a plain
# Find/Replace) so the corrected run keeps its original
# xml:space="preserve" attribute on <w:t>.
$typoIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "includingg migrgation") {
        $typoIndex = $i
        break
    }
}
$typoPara = $d.Paragraphs.Item($typoIndex)
$fixedXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1006"/><w:ilvl w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Introduction of a new backend payment system, including migration of the entire customer base. I was able to apply my previous experience from similar projects.</w:t></w:r></w:p>'
$typoPara.Range.InsertXML($fixedXml)

# --- 2. Insert the new "Ethos" bullet before the "Technical" heading --
# Locate the paragraph whose text is exactly "Technical" and which uses
# the Heading 2 style (the section heading carrying bookmark "technical"),
# and remember its 1-based index in the Paragraphs collection.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Technical" -and $p.Style.NameLocal -eq "Heading 2") {
        $targetIndex = $i
        break
    }
}

# Create a fresh, empty paragraph immediately before that heading - this
# shifts "Technical" (and everything after it) down by one paragraph
# without disturbing any of its own formatting/bookmarks. (Note: after
# this call the original paragraph object re-seats onto the new, empty
# paragraph rather than following "Technical", so re-fetch by index
# instead of reusing the object.)
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphBefore() | Out-Null

# The new, still-empty paragraph is now at $targetIndex (the "Technical"
# heading was pushed down to $targetIndex + 1). Fill it in (text +
# "Compact" style + the same bulleted-list numbering, numId 1002, used by
# the rest of the Summary bullets) via a raw OOXML fragment inserted into
# that paragraph's own (empty) range.
$newPara = $d.Paragraphs.Item($targetIndex)
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1002"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Ethos: Each line of code has an ongoing cost, so write as little as possible, and keep everything else tidy. Tests are first-class code.</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newXml)
